# Update "PERIOD TO EXPIRE" (column H) and "LAST UPDATE" (column I)
# for the Training Dashboard sheet to reflect progress as of 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 31; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $iCell = $ws.Cells.Item($row, 9)   # column I

    $currentPeriod = $hCell.Value()
    $hCell.Value = $currentPeriod - 1

    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
}
